# Apply the commit's data update: refresh "想去人数" (F) and "最低票价" (G)
# figures on the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4891
$ws1.Range("G3").Value  = 75
$ws1.Range("F5").Value  = 177
$ws1.Range("F6").Value  = 133
$ws1.Range("G6").Value  = 65
$ws1.Range("G7").Value  = 50
$ws1.Range("F10").Value = 246
$ws1.Range("F11").Value = 1255
$ws1.Range("F19").Value = 4224
$ws1.Range("F20").Value = 6524
$ws1.Range("F26").Value = 4035
$ws1.Range("F27").Value = 418
$ws1.Range("F35").Value = 334
$ws1.Range("F36").Value = 387
$ws1.Range("F37").Value = 199
$ws1.Range("F39").Value = 1584
$ws1.Range("F40").Value = 991
$ws1.Range("F41").Value = 54
$ws1.Range("F42").Value = 97
$ws1.Range("F48").Value = 604

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value  = 199

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 4891
$ws4.Range("G3").Value  = 75
$ws4.Range("F5").Value  = 177
$ws4.Range("F6").Value  = 133
$ws4.Range("G6").Value  = 65
$ws4.Range("G7").Value  = 199
$ws4.Range("G8").Value  = 50
$ws4.Range("F11").Value = 246
$ws4.Range("F12").Value = 1255
$ws4.Range("F20").Value = 4224
$ws4.Range("F21").Value = 6524
$ws4.Range("F27").Value = 4035
$ws4.Range("F28").Value = 418
$ws4.Range("F36").Value = 334
$ws4.Range("F37").Value = 387
$ws4.Range("F38").Value = 199
$ws4.Range("F40").Value = 1584
$ws4.Range("F41").Value = 991
$ws4.Range("F42").Value = 54
$ws4.Range("F43").Value = 97
$ws4.Range("F49").Value = 604
